# SDforEditStudentCommandLogicAndModel.pptx - resolved discrepancies in
# Sequence Diagram for EditStudentCommand.
#
# - refreshed the auto-updating "datetimeFigureOut" footer field (slide
#   master + all 11 slide layouts) from 29-10-2018 to 10-11-2018
# - repositioned / resized several shapes that make up the sequence
#   diagram (activation bars, labels, connectors) on slide 1
# - added a new dashed return-message arrow at the bottom of the diagram

# Shape.Left/Top/Width/Height round-trip through a single-precision
# (float32) "points" value in this COM host, and the EMU the value gets
# serialized back to is truncated rather than rounded - so converting an
# odd EMU count straight to points can land 1 EMU short after the
# round-trip. A tiny (sub-EMU) nudge before the points conversion keeps
# the truncated result exactly on the intended EMU value without
# perturbing values that already round-trip cleanly.
$EMU_PER_PT = 914400.0 / 72.0
function EMU($emu) { return ($emu + 0.5) / $EMU_PER_PT }

$p = $ppt.ActivePresentation

function Set-DatePlaceholderText($shapes, $newText) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            $tr = $shp.TextFrame.TextRange
            $len = $tr.Length
            if ($len -gt 0) {
                $chars = $tr.Characters(1, $len)
                $chars.Text = $newText
            }
        }
    }
}

# ---------------------------------------------------------------------
# 1) Refresh the datetimeFigureOut field everywhere it lives: the slide
#    master and every custom (slide) layout.
# ---------------------------------------------------------------------
$sm = $p.SlideMaster
Set-DatePlaceholderText $sm.Shapes "10-11-2018"
for ($li = 1; $li -le $sm.CustomLayouts.Count; $li++) {
    $cl = $sm.CustomLayouts.Item($li)
    Set-DatePlaceholderText $cl.Shapes "10-11-2018"
}

# ---------------------------------------------------------------------
# Slide 1 shape tweaks
# ---------------------------------------------------------------------
$s = $p.Slides.Item(1)

function Get-ShapeById($shapes, $id) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Id -eq $id) {
            return $shp
        }
    }
    return $null
}

# 2) "Rectangle 9" (id 10) - taller activation bar
$shp = Get-ShapeById $s.Shapes 10
$shp.Height = EMU(1577198)

# 3) "TextBox 16" (id 17, text "execute()") - moved down/right
$shp = Get-ShapeById $s.Shapes 17
$shp.Left = EMU(4322793)
$shp.Top = EMU(3790794)

# 4) "Straight Arrow Connector 18" (id 19) - moved down
$shp = Get-ShapeById $s.Shapes 19
$shp.Top = EMU(3235890)

# 5) "Straight Connector 38" (id 39) - taller activation bar connector
$shp = Get-ShapeById $s.Shapes 39
$shp.Width = EMU(3732)
$shp.Height = EMU(1451152)

# 6) "Rectangle 39" (id 40) - moved up 1 EMU, taller
$shp = Get-ShapeById $s.Shapes 40
$shp.Top = EMU(2212109)
$shp.Height = EMU(829537)

# 7) "Straight Arrow Connector 40" (id 41) - moved down/left
$shp = Get-ShapeById $s.Shapes 41
$shp.Left = EMU(3658955)
$shp.Top = EMU(3038028)

# 8) "Graphic 2" (id 3, the close icon picture) - moved down/right
$shp = Get-ShapeById $s.Shapes 3
$shp.Left = EMU(5217067)
$shp.Top = EMU(3358461)

# 9) "Connector: Elbow 34" (id 35) - resized/repositioned + adjustment
$shp = Get-ShapeById $s.Shapes 35
$shp.Left = EMU(1805294)
$shp.Top = EMU(2996508)
$shp.Width = EMU(6237477)
$shp.Height = EMU(1038199)
$shp.Adjustments.Item(1) = 0.92656

# ---------------------------------------------------------------------
# 10) New dashed return-message arrow at the bottom of the diagram,
#     styled like its sibling "Straight Arrow Connector" shapes.
# ---------------------------------------------------------------------
$template = Get-ShapeById $s.Shapes 73
$newShape = $template.Duplicate().Item(1)
$newShape.Name = "Straight Arrow Connector 60"
$newShape.Left = EMU(5429095)
$newShape.Top = EMU(2874776)
$newShape.Width = EMU(2698130)
$newShape.Height = EMU(0)
